# New crime data collected - weekly CompStat update (105th Precinct)
# Updates the report header (volume number + week-covering dates) and the
# crime-complaint statistics table (rows 15-30) to the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text: "Volume 32   Number  14" -> "...  15"
# ---------------------------------------------------------------------------
$ws.Range("A8").Characters(21, 2).Text = "15"

# ---------------------------------------------------------------------------
# Header text: "Report Covering the Week  3/31/2025  Through  4/6/2025"
#           -> "Report Covering the Week  4/7/2025  Through  4/13/2025"
# (Replace the first date before the second, recomputing the second date's
#  offset since its length changes after the first substitution.)
# ---------------------------------------------------------------------------
$ws.Range("C9").Characters(27, 9).Text = "4/7/2025"
$ws.Range("C9").Characters(46, 8).Text = "4/13/2025"

# ---------------------------------------------------------------------------
# Helper-ish inline approach: for cells that change numeric<->text "shape"
# we first clone formatting+value from a stable template cell of the
# desired shape (so the style index lines up with what Excel would reuse),
# then overwrite the value when the template's own value isn't the target.
# ---------------------------------------------------------------------------

# Row 15 (Murder): 28-day 2024 count 2 -> 1
$ws.Range("G15").Value = 1

# Row 16 (Robbery)
# C16: was text "0" -> becomes numeric 2 (adopt numeric style from F16 first)
$ws.Range("F16").Copy($ws.Range("C16"))
$ws.Range("C16").Value = 2
# D16: was numeric 1 -> becomes text "0" (clone style+value from template C27)
$ws.Range("C27").Copy($ws.Range("D16"))
# E16: was numeric -100 -> becomes text "***.*" (clone from template E27)
$ws.Range("E27").Copy($ws.Range("E16"))
$ws.Range("F16").Value = 4
$ws.Range("H16").Value = 100
$ws.Range("I16").Value = 25
$ws.Range("K16").Value = 78.571428571428
$ws.Range("L16").Value = 8.695652173913
$ws.Range("M16").Value = -71.590909090909
$ws.Range("N16").Value = -90.842490842490

# Row 17 (Fel. Assault)
$ws.Range("C17").Value = 3
# D17: was text "0" -> becomes numeric 8 (adopt numeric style from F17 first)
$ws.Range("F17").Copy($ws.Range("D17"))
$ws.Range("D17").Value = 8
# E17: was text "***.*" -> becomes numeric -62.5 (adopt numeric style from H17 first)
$ws.Range("H17").Copy($ws.Range("E17"))
$ws.Range("E17").Value = -62.5
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = 12.5
$ws.Range("I17").Value = 59
$ws.Range("J17").Value = 80
$ws.Range("K17").Value = -26.25
$ws.Range("L17").Value = -1.666666666666
$ws.Range("M17").Value = -23.376623376623
$ws.Range("N17").Value = -35.869565217391

# Row 18 (Burglary)
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 50
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = 15.384615384615
$ws.Range("I18").Value = 48
$ws.Range("J18").Value = 38
$ws.Range("K18").Value = 26.315789473684
$ws.Range("L18").Value = 17.073170731707
$ws.Range("M18").Value = -54.285714285714
$ws.Range("N18").Value = -88.652482269503

# Row 19 (Gr. Larceny)
$ws.Range("C19").Value = 4
$ws.Range("E19").Value = -55.555555555555
$ws.Range("F19").Value = 30
$ws.Range("G19").Value = 35
$ws.Range("H19").Value = -14.285714285714
$ws.Range("I19").Value = 105
$ws.Range("J19").Value = 118
$ws.Range("K19").Value = -11.016949152542
$ws.Range("L19").Value = -0.943396226415
$ws.Range("M19").Value = -19.230769230769
$ws.Range("N19").Value = -31.818181818181

# Row 20 (G.L.A.)
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 22
$ws.Range("G20").Value = 18
$ws.Range("H20").Value = 22.222222222222
$ws.Range("I20").Value = 53
$ws.Range("J20").Value = 59
$ws.Range("K20").Value = -10.169491525423
$ws.Range("L20").Value = 10.416666666666
$ws.Range("M20").Value = -50.467289719626
$ws.Range("N20").Value = -94.569672131147

# Row 21 (TOTAL)
$ws.Range("C21").Value = 15
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = -31.818181818181
$ws.Range("F21").Value = 89
$ws.Range("G21").Value = 85
$ws.Range("H21").Value = 4.705882352941
$ws.Range("I21").Value = 294
$ws.Range("J21").Value = 315
$ws.Range("K21").Value = -6.666666666666
$ws.Range("L21").Value = 4.626334519572
$ws.Range("M21").Value = -43.133462282398
$ws.Range("N21").Value = -84.806201550387

# Row 24 (Petit Larceny)
$ws.Range("C24").Value = 4
$ws.Range("D24").Value = 10
$ws.Range("E24").Value = -60
$ws.Range("F24").Value = 29
$ws.Range("G24").Value = 42
$ws.Range("H24").Value = -30.952380952381
$ws.Range("I24").Value = 148
$ws.Range("J24").Value = 197
$ws.Range("K24").Value = -24.873096446700
$ws.Range("L24").Value = -33.928571428571
$ws.Range("M24").Value = -34.222222222222

# Row 25 (Retail Theft)
$ws.Range("F25").Value = 6
$ws.Range("G25").Value = 7
$ws.Range("H25").Value = -14.285714285714
$ws.Range("I25").Value = 33
$ws.Range("J25").Value = 44
$ws.Range("K25").Value = -25
$ws.Range("L25").Value = -10.810810810810

# Row 26 (Misd. Assault)
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 140
$ws.Range("F26").Value = 43
$ws.Range("H26").Value = 95.454545454545
$ws.Range("I26").Value = 124
$ws.Range("J26").Value = 92
$ws.Range("K26").Value = 34.782608695652
$ws.Range("L26").Value = 42.528735632183
$ws.Range("M26").Value = -29.943502824858

# Row 27 (UCR Rape*)
$ws.Range("G27").Value = 1

# Row 28 (Other Sex Crimes)
# C28: was text "0" -> becomes numeric 1 (adopt numeric style from F28 first)
$ws.Range("F28").Copy($ws.Range("C28"))
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 1
$ws.Range("H28").Value = -75
$ws.Range("I28").Value = 11
$ws.Range("J28").Value = 14
$ws.Range("K28").Value = -21.428571428571
$ws.Range("L28").Value = 266.666666666667

# Row 29 (Shooting Vic.)
# G29: was numeric 1 -> becomes text "0" (clone style+value from template F29)
$ws.Range("F29").Copy($ws.Range("G29"))
# H29: was numeric -100 -> becomes text "***.*" (clone from template E29)
$ws.Range("E29").Copy($ws.Range("H29"))

# Row 30 (Shooting Inc.)
$ws.Range("F30").Copy($ws.Range("G30"))
$ws.Range("E30").Copy($ws.Range("H30"))
